$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: Status text "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: Status text + Latest Handback DateTime refresh + clear Error Detail ---
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-10-18 04:03:27"
$zhcn.Range("K3").Value = "2016-10-18 04:03:27"
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# --- de-de sheet: Status text + Latest Handback DateTime refresh + clear Error Detail ---
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-10-18 04:03:50"
$dede.Range("K3").Value = "2016-10-18 04:03:50"
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""
